$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.385.02"
$ws.Range("E2").Value = "  +3.57%  "

$ws.Range("D3").Value = "3.152.47"
$ws.Range("E3").Value = "  +2.97%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "397.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.548"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.75%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.84%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.70%  "

$ws.Range("E11").Value = "  +1.26%  "

$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("D13").Value = "3.648.42"
$ws.Range("E13").Value = "  +2.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.50%  "

$ws.Range("E16").Value = "  +8.87%  "

$ws.Range("D17").Value = "3.150.37"
$ws.Range("E17").Value = "  +2.93%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.25%  "

$ws.Range("D19").Value = "53.329.89"
$ws.Range("E19").Value = "  +3.28%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.07%  "

$ws.Range("D22").Value = "0.0₃0974"
$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "270.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("E25").Value = "  +2.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.96%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.77%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("E31").Value = "  +1.99%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "37.33"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0499"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.80%  "

$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.67"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.01%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.74%  "

$ws.Range("E41").Value = "  -0.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.44%  "

$ws.Range("E43").Value = "  +1.64%  "

$ws.Range("E44").Value = "  +4.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.75%  "

$ws.Range("E47").Value = "  -1.34%  "

$ws.Range("D48").Value = "2.089.84"
$ws.Range("E48").Value = "  +2.74%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +21.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.70%  "
